$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 (Varta Ultra Lithium AA 4er Bli), shifting subsequent rows up
$ws.Rows.Item(3).Delete()

# Update the timestamp column (O) for the remaining data rows (2-29) to the new crawl time
$ws.Range("O2:O29").Value = "2022-07-24 20:58:53"
